$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

# Column E (Sprint 2 status) - copy the "PASS" cell E10 (value + style) down
# into E11:E18, matching how the author pasted the status into each row.
for ($r = 11; $r -le 18; $r++) {
    $ws.Range("E10").Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
    $ws.Range("E10").Copy()
    $ws.Cells.Item($r, 5).PasteSpecial()
}

# Column F (Sprint 2 comments) - type the tester initials/date directly;
# this keeps each cell's pre-existing style, only adding the text.
for ($r = 11; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = "JE; 4/1/2018"
}
